# Update the single-column results table on the Renaissance / ZGC
# movie-lens (heap 8G) benchmark doc.
#
# Rows 1-12 (1-indexed) hold individual summary values that were
# recomputed; rows 44-46 used to carry a full tab-separated line of
# per-iteration numbers but now only keep the first (headline) value,
# matching the single value already used for every other result row
# in the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "819"
$t.Cell(5, 1).Range.Text  = "0.00002"
$t.Cell(6, 1).Range.Text  = "0.00039"
$t.Cell(7, 1).Range.Text  = "0.00005"
$t.Cell(8, 1).Range.Text  = "0.00002"
$t.Cell(9, 1).Range.Text  = "0.00004"
$t.Cell(10, 1).Range.Text = "0.00005"
$t.Cell(11, 1).Range.Text = "0.00006"
$t.Cell(12, 1).Range.Text = "0.04220"

$t.Cell(44, 1).Range.Text = "100"
$t.Cell(45, 1).Range.Text = "0.04"
$t.Cell(46, 1).Range.Text = "3861"
